# Rewrites the opening paragraph of the document:
#   "This is a test document. With a bgu. Oh noes"
# becomes:
#   "This is a test document. Without a bug. Oh boy!"
#
# The target OOXML splits the new sentence across four separate <w:r> runs
# ("Th" / "is is a test document. Without a b" / "u" / "g. Oh boy!") and
# drops the <w:proofErr> spell-check markers that bracketed the old "bgu".
# The trailing <w:bookmarkStart/bookmarkEnd w:name="_GoBack"/> pair must
# remain, now sitting after the new text instead of after the old text.

$d = $word.ActiveDocument

# The four pieces of the rebuilt sentence, and their cumulative offsets
# from the start of the paragraph (0-based, computed up front so we never
# depend on a Range object "growing" after an Insert call).
$piece1 = "Th"
$piece2 = "is is a test document. Without a b"
$piece3 = "u"
$piece4 = "g. Oh boy!"

$off0 = 0
$off1 = $off0 + $piece1.Length
$off2 = $off1 + $piece2.Length
$off3 = $off2 + $piece3.Length
$off4 = $off3 + $piece4.Length

# Locate the paragraph's original text span (everything up to, but not
# including, the paragraph mark).
$para = $d.Paragraphs.Item(1)
$oldTextEnd = $para.Range.End - 1
$oldTextRange = $d.Range(0, $oldTextEnd)

# Remove the old sentence (and, with it, the spellcheck proofErr markers
# that live between its runs) while leaving the "_GoBack" bookmark alone.
$oldTextRange.Delete()

# Rebuild the sentence as four discrete runs by inserting each piece at
# the boundary of the previous one. Using InsertAfter on a freshly
# constructed Range (rather than the bookmark's own Range) keeps each
# inserted piece as its own run instead of merging into its neighbor.
$d.Range($off0, $off0).InsertAfter($piece1)
$d.Range($off1, $off1).InsertAfter($piece2)
$d.Range($off2, $off2).InsertAfter($piece3)
$d.Range($off3, $off3).InsertAfter($piece4)

# A temporary sacrificial character is appended after the new text so
# that, when the bookmark is re-created, its position is not the very
# last character of the story -- re-adding "_GoBack" exactly at
# end-of-text mis-places it at the start of the paragraph instead.
$d.Range($off4, $off4).InsertAfter("Z")

$bookmark = $d.Bookmarks.Item("_GoBack")
$bookmark.Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($off4, $off4))

# Drop the sacrificial character now that the bookmark sits correctly
# right after the new text.
$d.Range($off4, $off4 + 1).Delete()
